$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the numeric-looking columns (Price, Volume(1h), Hora)
# so Excel keeps them as text strings instead of auto-converting to numbers/percentages.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "323.59"
$ws.Range("E2").Value = "1.12%"
$ws.Range("G2").Value = "2"

# Row 3
$ws.Range("D3").Value = "39.42"
$ws.Range("E3").Value = "-0.77%"
$ws.Range("G3").Value = "2"

# Row 4
$ws.Range("D4").Value = "5.959"
$ws.Range("E4").Value = "13.66%"
$ws.Range("G4").Value = "2"

# Row 5
$ws.Range("D5").Value = "0.08016"
$ws.Range("E5").Value = "-1.18%"
$ws.Range("G5").Value = "2"

# Row 6
$ws.Range("D6").Value = "4.590"
$ws.Range("E6").Value = "1.75%"
$ws.Range("G6").Value = "2"

# Row 7
$ws.Range("D7").Value = "8.654"
$ws.Range("E7").Value = "0.76%"
$ws.Range("G7").Value = "2"

# Row 8
$ws.Range("D8").Value = "1.913"
$ws.Range("E8").Value = "0.29%"
$ws.Range("G8").Value = "2"

# Row 9
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9354"
$ws.Range("E9").Value = "0.08%"
$ws.Range("G9").Value = "2"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1257"
$ws.Range("E10").Value = "-2.21%"
$ws.Range("G10").Value = "2"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1972"
$ws.Range("E11").Value = "0.53%"
$ws.Range("G11").Value = "2"

# Row 12
$ws.Range("B12").Value = "MCDex"
$ws.Range("C12").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D12").Value = "8.795"
$ws.Range("E12").Value = "30.10%"
$ws.Range("G12").Value = "2"

# Row 13
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D13").Value = "0.09242"
$ws.Range("E13").Value = "-0.31%"
$ws.Range("G13").Value = "2"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03421"
$ws.Range("E14").Value = "0.58%"
$ws.Range("G14").Value = "2"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09601"
$ws.Range("E15").Value = "0.87%"
$ws.Range("G15").Value = "2"

# Row 16
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001289"
$ws.Range("E16").Value = "-7.27%"
$ws.Range("G16").Value = "2"

# Row 17
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "0.006325"
$ws.Range("E17").Value = "7.64%"
$ws.Range("G17").Value = "2"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "3.338"
$ws.Range("E18").Value = "-0.88%"
$ws.Range("G18").Value = "2"

# Row 19
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.941"
$ws.Range("E19").Value = "-1.96%"
$ws.Range("G19").Value = "2"

# Row 20
$ws.Range("D20").Value = "0.3536"
$ws.Range("E20").Value = "0.05%"
$ws.Range("G20").Value = "2"

# Row 21
$ws.Range("D21").Value = "0.1400"
$ws.Range("E21").Value = "5.43%"
$ws.Range("G21").Value = "2"

# Row 22
$ws.Range("D22").Value = "0.2411"
$ws.Range("E22").Value = "4.31%"
$ws.Range("G22").Value = "2"

# Row 23
$ws.Range("D23").Value = "0.04467"
$ws.Range("E23").Value = "0.66%"
$ws.Range("G23").Value = "2"

# Row 24
$ws.Range("D24").Value = "0.001263"
$ws.Range("E24").Value = "3.30%"
$ws.Range("G24").Value = "2"

# Row 25
$ws.Range("D25").Value = "0.004367"
$ws.Range("E25").Value = "0.19%"
$ws.Range("G25").Value = "2"

# Row 26
$ws.Range("D26").Value = "0.0001141"
$ws.Range("E26").Value = "-11.68%"
$ws.Range("G26").Value = "2"

# Row 27
$ws.Range("E27").Value = "0.07%"
$ws.Range("G27").Value = "2"

# Row 28
$ws.Range("G28").Value = "2"

# Row 29
$ws.Range("G29").Value = "2"

# Row 30
$ws.Range("G30").Value = "2"

# Row 31
$ws.Range("G31").Value = "2"

# Row 32
$ws.Range("G32").Value = "2"

# Row 33
$ws.Range("G33").Value = "2"

# Row 34
$ws.Range("G34").Value = "2"

# Row 35
$ws.Range("G35").Value = "2"

# Row 36
$ws.Range("G36").Value = "2"

# Row 37
$ws.Range("G37").Value = "2"

# Row 38
$ws.Range("G38").Value = "2"

# Row 39
$ws.Range("D39").Value = "0.02431"
$ws.Range("E39").Value = "-0.19%"
$ws.Range("G39").Value = "2"

# Row 40
$ws.Range("D40").Value = "0.05207"
$ws.Range("E40").Value = "-0.05%"
$ws.Range("G40").Value = "2"

# Row 41
$ws.Range("D41").Value = "0.007444"
$ws.Range("E41").Value = "-2.68%"
$ws.Range("G41").Value = "2"

# Row 42
$ws.Range("D42").Value = "0.1409"
$ws.Range("E42").Value = "-1.66%"
$ws.Range("G42").Value = "2"

# Row 43
$ws.Range("D43").Value = "0.008893"
$ws.Range("E43").Value = "3.01%"
$ws.Range("G43").Value = "2"

# Row 44
$ws.Range("D44").Value = "0.002095"
$ws.Range("E44").Value = "-0.88%"
$ws.Range("G44").Value = "2"

# Row 45
$ws.Range("D45").Value = "0.01126"
$ws.Range("E45").Value = "25.34%"
$ws.Range("G45").Value = "2"

# Row 46
$ws.Range("D46").Value = "0.00006728"
$ws.Range("E46").Value = "2.27%"
$ws.Range("G46").Value = "2"

# Row 47
$ws.Range("E47").Value = "0.09%"
$ws.Range("G47").Value = "2"

# Row 48
$ws.Range("D48").Value = "0.003004"
$ws.Range("E48").Value = "5.03%"
$ws.Range("G48").Value = "2"

# Row 49
$ws.Range("E49").Value = "-42.85%"
$ws.Range("G49").Value = "2"

# Row 50
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.09%"
$ws.Range("G50").Value = "2"

# Row 51
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.09%"
$ws.Range("G51").Value = "2"
